# Change modeling of rest of the world:
# - there can only be one rest-of-world region in the model aggregation
# - trade with rest of the world is shifted into the overall trade block;
#   now there is a distinction between total import (including from rest
#   of the world) and import from modeled regions.
#
# Concretely, on the "elasTRADE" sheet:
#   - old column B "elasIMP" is dropped (replaced by the shifted headers
#     below), old column C "elasIU_DM" becomes column B, old column D
#     "elasFU_DM" becomes column C
#   - a new column D "elasIMP_ROW" (import from rest-of-world-excluded /
#     modeled regions) is added
#   - a new column E "elasTRD" (overall trade block) is added, filled
#     with the same constant (5) as the other trade columns
# The elasTRADE sheet also becomes the active/selected sheet instead of TFP.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("elasTRADE")

# Shift the existing header labels left by one column.
$ws.Range("B1").Value = "elasIU_DM"
$ws.Range("C1").Value = "elasFU_DM"

# New columns appended at the end - add the "elasTRD" header first so it
# lands earlier in the shared-string table than "elasIMP_ROW", matching
# how the strings were originally authored.
$ws.Range("E1").Value = "elasTRD"
$ws.Range("D1").Value = "elasIMP_ROW"

# Fill the new "elasTRD" column (E) with the same constant value (5)
# used throughout the rest of the trade-elasticity table.
for ($r = 2; $r -le 36; $r++) {
  $ws.Cells.Item($r, 5).Value = 5
}

# elasTRADE becomes the active sheet (was TFP before).
$ws.Activate()
